$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.372.90'
$ws.Range("E2").Value = '  -2.84%  '
$ws.Range("D3").Value = '2.220.73'
$ws.Range("E3").Value = '  -6.15%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '296.47'
$ws.Range("E5").Value = '  -4.60%  '
$ws.Range("D6").Value = '83.46'
$ws.Range("E6").Value = '  -2.96%  '
$ws.Range("D7").Value = '0.511'
$ws.Range("E7").Value = '  -3.15%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").Value = '0.469'
$ws.Range("E9").Value = '  -4.35%  '
$ws.Range("D10").Value = '0.0776'
$ws.Range("E10").Value = '  -7.58%  '
$ws.Range("D11").Value = '29.21'
$ws.Range("E11").Value = '  -2.95%  '
$ws.Range("D12").Formula = "'47.70"
$ws.Range("E12").Value = '  -9.57%  '
$ws.Range("E13").Value = '  -1.96%  '
$ws.Range("D14").Value = '2.558.51'
$ws.Range("E14").Value = '  -6.37%  '
$ws.Range("D15").Value = '6.29'
$ws.Range("E15").Value = '  -3.76%  '
$ws.Range("D16").Value = '14.17'
$ws.Range("E16").Value = '  -5.65%  '
$ws.Range("D17").Value = '2.200.41'
$ws.Range("E17").Value = '  -7.43%  '
$ws.Range("D18").Value = '0.715'
$ws.Range("E18").Value = '  -5.23%  '
$ws.Range("D19").Value = '39.287.70'
$ws.Range("E19").Value = '  -2.99%  '
$ws.Range("D20").Value = '0.0₃0874'
$ws.Range("E20").Value = '  -3.97%  '
$ws.Range("D21").Value = '5.72'
$ws.Range("E21").Value = '  -6.40%  '
$ws.Range("D22").Value = '64.92'
$ws.Range("E22").Value = '  -5.19%  '
$ws.Range("D23").Value = '10.44'
$ws.Range("E23").Value = '  -2.50%  '
$ws.Range("D24").Value = '228.43'
$ws.Range("E24").Value = '  -2.75%  '
$ws.Range("E25").Value = '  +0.13%  '
$ws.Range("D26").Value = '2.41'
$ws.Range("E26").Value = '  -6.46%  '
$ws.Range("D27").Formula = "'1.80"
$ws.Range("E27").Value = '  +0.08%  '
$ws.Range("D28").Value = '22.59'
$ws.Range("E28").Value = '  -4.94%  '
$ws.Range("E29").Value = '  -2.61%  '
$ws.Range("D30").Value = '9.13'
$ws.Range("E30").Value = '  -1.07%  '
$ws.Range("D31").Value = '149.85'
$ws.Range("E31").Value = '  -2.72%  '
$ws.Range("D32").Value = '31.94'
$ws.Range("E32").Value = '  -6.52%  '
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  -0.27%  '
$ws.Range("D34").Value = '4.86'
$ws.Range("E34").Value = '  -6.16%  '
$ws.Range("D35").Value = '0.0695'
$ws.Range("E35").Value = '  -4.16%  '
$ws.Range("D36").Value = '2.33'
$ws.Range("E36").Value = '  -3.35%  '
$ws.Range("E37").Value = '  -3.41%  '
$ws.Range("B38").Value = 'Celestia'
$ws.Range("C38").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D38").Value = '15.32'
$ws.Range("E38").Value = '  -4.69%  '
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").Value = '0.0963'
$ws.Range("E39").Value = '  -3.95%  '
$ws.Range("D40").Value = '2.64'
$ws.Range("E40").Value = '  -3.62%  '
$ws.Range("D41").Value = '1.63'
$ws.Range("E41").Value = '  -4.24%  '
$ws.Range("D42").Value = '3.67'
$ws.Range("E42").Value = '  -4.21%  '
$ws.Range("D43").Value = '1.910.23'
$ws.Range("E43").Value = '  -2.38%  '
$ws.Range("D44").Value = '0.0259'
$ws.Range("E44").Value = '  -3.33%  '
$ws.Range("E45").Value = '  -15.81%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '16.02'
$ws.Range("E46").Value = '  -8.92%  '
$ws.Range("D47").Value = '2.63'
$ws.Range("E47").Value = '  -1.87%  '
$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").Value = '8.91'
$ws.Range("E48").Value = '  -4.41%  '
$ws.Range("D49").Value = '2.415.24'
$ws.Range("E49").Value = '  -7.07%  '
$ws.Range("D50").Value = '70.68'
$ws.Range("E50").Value = '  -1.52%  '
$ws.Range("D51").Value = '87.22'
$ws.Range("E51").Value = '  -6.06%  '
